$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update F column "想去人数" counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1624
$ws1.Range("F6").Value = 23
$ws1.Range("F7").Value = 415
$ws1.Range("F10").Value = 485

# Sheet "全部类型" (sheet4): update F column "想去人数" counts
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 387
$ws4.Range("F3").Value = 115
$ws4.Range("F4").Value = 1624
$ws4.Range("F5").Value = 17
$ws4.Range("F7").Value = 415
$ws4.Range("F9").Value = 63
$ws4.Range("F10").Value = 485

$wb.Save()
